# [지현] buffitem fielditem 이름 맞춤
# Rename FieldItemData.itemCode values to match BuffItemData.itemCode values
# (append "Buff" suffix so field items line up with their buff-item counterparts).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("FieldItemData")

$ws.Range("A6").Value  = "MinusDiceBuff"
$ws.Range("A7").Value  = "DoubleDiceBuff"
$ws.Range("A8").Value  = "HalfDiceBuff"
$ws.Range("A9").Value  = "DrunkBuff"
$ws.Range("A10").Value = "OddBuff"
$ws.Range("A11").Value = "EvenBuff"

# Restore the selection state recorded in the saved workbook for BuffItemData
# (without leaving it as the active/selected sheet).
$wsBuff = $wb.Worksheets.Item("BuffItemData")
$wsBuff.Select()
$wsBuff.Range("A3:A8").Select()
$wsBuff.Range("A8").Activate()

# FieldItemData stays the active sheet/tab, with its own selection restored.
$ws.Select()
$ws.Range("B18").Select()
